$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "43.700.25"
$ws.Range("E2").Value = "  +2.47%  "

$ws.Range("D3").Value = "2.234.99"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.03"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "79.19"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.03%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.13"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0925"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.11"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.31%  "

$ws.Range("E13").Value = "  +1.46%  "

$ws.Range("D14").Value = "2.580.55"
$ws.Range("E14").Value = "  +1.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.73"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.22%  "

$ws.Range("D16").Value = "2.250.11"
$ws.Range("E16").Value = "  +1.65%  "

$ws.Range("E17").Value = "  +2.09%  "

$ws.Range("D18").Value = "43.599.43"
$ws.Range("E18").Value = "  +2.44%  "

$ws.Range("E19").Value = "  +2.53%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.54"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.06"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.94%  "

$ws.Range("E22").Value = "  +7.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.91"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.45"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "42.55"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.84"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.37"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.77%  "

$ws.Range("E29").Value = "  +1.37%  "

$ws.Range("E30").Value = "  -1.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.11"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.62"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0875"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.04%  "

$ws.Range("E34").Value = "  +2.20%  "

$ws.Range("E35").Value = "  +1.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0365"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.50"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.108"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.16"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.90"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +20.85%  "

$ws.Range("E41").Value = "  +2.56%  "

$ws.Range("E42").Value = "  +1.58%  "

$ws.Range("E43").Value = "  +5.50%  "

$ws.Range("E44").Value = "  +2.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.67"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.53"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.473"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0988"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.12"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.70%  "

$ws.Range("E50").Value = "  +2.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.52"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +27.38%  "

